$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '34.632.34'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.96%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.817.27'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.44%  '

$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '228.66'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.09%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.560'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.77%  '

$ws.Range("E7").Value = '  +0.13%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '34.90'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +7.87%  '

$ws.Range("E9").Value = '  +1.68%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0695'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.78%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0951'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.27%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.078.24'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.42%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.36'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.82%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.812.30'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.23%  '

$ws.Range("E15").Value = '  +2.63%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '34.662.34'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.12%  '

$ws.Range("E17").Value = '  +3.13%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '69.23'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.68%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '247.51'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.31%  '

$ws.Range("E20").Value = '  -0.15%  '

$ws.Range("E21").Value = '  +5.47%  '

$ws.Range("E22").Value = '  +0.17%  '

$ws.Range("E23").Value = '  +0.88%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '172.13'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +6.00%  '

$ws.Range("E25").Value = '  +1.22%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.46'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.80%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.78'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.52%  '

$ws.Range("E28").Value = '  +1.53%  '

$ws.Range("E29").Value = '  -0.14%  '

$ws.Range("E30").Value = '  +4.02%  '

$ws.Range("E31").Value = '  +2.13%  '

$ws.Range("E32").Value = '  +2.41%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.25'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.40%  '

$ws.Range("E34").Value = '  +2.75%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.62'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.40%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.420.35'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.60%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.679'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.72%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.07'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.40%  '

$ws.Range("B39").Value = 'Aave'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '86.37'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.89%  '

$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0192'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.92%  '

$ws.Range("E41").Value = '  +4.43%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.965'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.38%  '

$ws.Range("E43").Value = '  +0.47%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.87'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.97%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0524'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.80%  '

$ws.Range("E46").Value = '  +3.02%  '

$ws.Range("E47").Value = '  +0.88%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.979.41'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.73%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '106.29'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.68%  '

$ws.Range("E50").Value = '  +0.58%  '

$ws.Range("E51").Value = '  +0.10%  '
